# Auto-generated edit script: updates Profit-table formula results (H/I/J/K/L/M/N columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets, per the source diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H5").Value = 136.75
$ws.Range("I5").Value = 136.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 136.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -21.75
$ws.Range("N5").ClearContents()
$ws.Range("H17").Value = 2021
$ws.Range("I17").Value = 1531.5
$ws.Range("J17").Value = 3000
$ws.Range("K17").Value = 4594.5
$ws.Range("L17").Value = 9000
$ws.Range("M17").Value = -4426.5
$ws.Range("N17").Value = -9336
$ws.Range("H19").Value = 805.4286
$ws.Range("I19").Value = 659.9
$ws.Range("K19").Value = 659.9
$ws.Range("M19").Value = -484.9
$ws.Range("H53").Value = 252.91667
$ws.Range("I53").Value = 399.66666
$ws.Range("J53").Value = 204
$ws.Range("K53").Value = 399.66666
$ws.Range("L53").Value = 204
$ws.Range("M53").Value = 237.33334
$ws.Range("N53").Value = -1478
$ws.Range("H95").Value = 83002.5
$ws.Range("J95").Value = 83002.5
$ws.Range("L95").Value = 83002.5
$ws.Range("N95").Value = -88494.5
$ws.Range("H97").Value = 2129.75
$ws.Range("J97").Value = 2129.75
$ws.Range("L97").Value = 6389.25
$ws.Range("N97").Value = -7381.25
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H112").Value = 2466
$ws.Range("I112").Value = 2466
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 7398
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -6290
$ws.Range("N112").ClearContents()
$ws.Range("H113").Value = 4318.9
$ws.Range("I113").Value = 3779.8
$ws.Range("K113").Value = 3779.8
$ws.Range("M113").Value = -525.8000000000002
$ws.Range("H132").Value = 17452.75
$ws.Range("I132").Value = 19641.846
$ws.Range("J132").Value = 7966.6665
$ws.Range("K132").Value = 58925.538
$ws.Range("L132").Value = 23899.9995
$ws.Range("M132").Value = -56395.538
$ws.Range("N132").Value = -28959.9995
$ws.Range("H137").Value = 14022.154
$ws.Range("I137").Value = 2912.7144
$ws.Range("J137").Value = 26983.166
$ws.Range("K137").Value = 8738.143199999999
$ws.Range("L137").Value = 80949.498
$ws.Range("M137").Value = -6188.143199999999
$ws.Range("N137").Value = -86049.498
$ws.Range("H138").Value = 2171.419
$ws.Range("I138").Value = 1243.2941
$ws.Range("J138").Value = 2448.228
$ws.Range("K138").Value = 3729.8823
$ws.Range("L138").Value = 7344.684
$ws.Range("M138").Value = 1410.1177
$ws.Range("N138").Value = -17624.684

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 213766.27
$ws.Range("I32").Value = 222542.83
$ws.Range("K32").Value = 222542.83
$ws.Range("M32").Value = -222255.83
$ws.Range("H45").Value = 3540.8572
$ws.Range("J45").Value = 3395
$ws.Range("L45").Value = 3395
$ws.Range("N45").Value = -4149
$ws.Range("H61").Value = 4538.4644
$ws.Range("I61").Value = 4281.773
$ws.Range("J61").Value = 5479.6665
$ws.Range("K61").Value = 4281.773
$ws.Range("L61").Value = 5479.6665
$ws.Range("M61").Value = -4069.773
$ws.Range("N61").Value = -5903.6665
$ws.Range("H122").Value = 2121.5625
$ws.Range("I122").Value = 1967.5714
$ws.Range("K122").Value = 5902.7142
$ws.Range("M122").Value = -3452.7142
$ws.Range("H132").Value = 659891.3
$ws.Range("I132").Value = 759521.5
$ws.Range("J132").Value = 2331.8
$ws.Range("K132").Value = 2278564.5
$ws.Range("L132").Value = 6995.400000000001
$ws.Range("M132").Value = -2276034.5
$ws.Range("N132").Value = -12055.4
$ws.Range("H133").Value = 52222.145
$ws.Range("J133").Value = 52222.145
$ws.Range("L133").Value = 52222.145
$ws.Range("N133").Value = -57282.145
$ws.Range("H134").Value = 60427.5
$ws.Range("J134").Value = 60427.5
$ws.Range("L134").Value = 60427.5
$ws.Range("N134").Value = -70567.5
$ws.Range("H136").Value = 4538.4644
$ws.Range("I136").Value = 4281.773
$ws.Range("J136").Value = 5479.6665
$ws.Range("K136").Value = 12845.319
$ws.Range("L136").Value = 16438.9995
$ws.Range("M136").Value = -10295.319
$ws.Range("N136").Value = -21538.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2176.5
$ws.Range("I86").Value = 2402
$ws.Range("K86").Value = 2402
$ws.Range("M86").Value = -1279
$ws.Range("H89").Value = 2176.5
$ws.Range("I89").Value = 2402
$ws.Range("K89").Value = 12010
$ws.Range("M89").Value = -6394
$ws.Range("H134").Value = 6039.3667
$ws.Range("I134").Value = 2699.0952
$ws.Range("K134").Value = 8097.285600000001
$ws.Range("M134").Value = -5562.285600000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 2014.5714
$ws.Range("I25").Value = 2014.5714
$ws.Range("K25").Value = 2014.5714
$ws.Range("M25").Value = -1840.5714
$ws.Range("H31").Value = 2471.4736
$ws.Range("I31").Value = 2419.9092
$ws.Range("J31").Value = 2542.375
$ws.Range("K31").Value = 2419.9092
$ws.Range("L31").Value = 2542.375
$ws.Range("M31").Value = -2124.9092
$ws.Range("N31").Value = -3132.375
$ws.Range("H34").Value = 2471.4736
$ws.Range("I34").Value = 2419.9092
$ws.Range("J34").Value = 2542.375
$ws.Range("K34").Value = 2419.9092
$ws.Range("L34").Value = 2542.375
$ws.Range("M34").Value = -2217.9092
$ws.Range("N34").Value = -2946.375
$ws.Range("H58").Value = 30003.25
$ws.Range("I58").Value = 6999
$ws.Range("J58").Value = 37671.332
$ws.Range("K58").Value = 6999
$ws.Range("L58").Value = 37671.332
$ws.Range("M58").Value = -6796
$ws.Range("N58").Value = -38077.332
$ws.Range("H107").Value = 697.7692
$ws.Range("I107").Value = 540.25
$ws.Range("K107").Value = 540.25
$ws.Range("M107").Value = 1379.75
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 10358.1
$ws.Range("I122").Value = 1690.3846
$ws.Range("J122").Value = 26455.285
$ws.Range("K122").Value = 5071.1538
$ws.Range("L122").Value = 79365.855
$ws.Range("M122").Value = -2621.1538
$ws.Range("N122").Value = -84265.855
$ws.Range("H127").Value = 90000
$ws.Range("J127").Value = 90000
$ws.Range("L127").Value = 90000
$ws.Range("N127").Value = -99920
$ws.Range("H132").Value = 4066.4
$ws.Range("I132").Value = 4006
$ws.Range("K132").Value = 12018
$ws.Range("M132").Value = -9488
$ws.Range("H134").Value = 3724.875
$ws.Range("I134").Value = 2559.8
$ws.Range("K134").Value = 7679.400000000001
$ws.Range("M134").Value = -5144.400000000001
$ws.Range("H136").Value = 30003.25
$ws.Range("I136").Value = 6999
$ws.Range("J136").Value = 37671.332
$ws.Range("K136").Value = 20997
$ws.Range("L136").Value = 113013.996
$ws.Range("M136").Value = -18447
$ws.Range("N136").Value = -118113.996

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 248.27272
$ws.Range("J26").Value = 424.83334
$ws.Range("L26").Value = 1274.50002
$ws.Range("N26").Value = -1850.50002
$ws.Range("H36").Value = 505.83334
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H62").Value = 1418.75
$ws.Range("I62").Value = 891.6667
$ws.Range("K62").Value = 2675.0001
$ws.Range("M62").Value = -1989.0001
$ws.Range("H63").Value = 2200
$ws.Range("I63").Value = 2200
$ws.Range("K63").Value = 6600
$ws.Range("M63").Value = -5851
$ws.Range("H64").Value = 5332.1665
$ws.Range("I64").Value = 4249.5
$ws.Range("J64").Value = 7497.5
$ws.Range("K64").Value = 12748.5
$ws.Range("L64").Value = 22492.5
$ws.Range("M64").Value = -12478.5
$ws.Range("N64").Value = -23032.5
$ws.Range("H65").Value = 1418.75
$ws.Range("I65").Value = 891.6667
$ws.Range("K65").Value = 8025.0003
$ws.Range("M65").Value = -4593.0003
$ws.Range("H66").Value = 2200
$ws.Range("I66").Value = 2200
$ws.Range("K66").Value = 19800
$ws.Range("M66").Value = -16056
$ws.Range("H67").Value = 5332.1665
$ws.Range("I67").Value = 4249.5
$ws.Range("J67").Value = 7497.5
$ws.Range("K67").Value = 12748.5
$ws.Range("L67").Value = 22492.5
$ws.Range("M67").Value = -11812.5
$ws.Range("N67").Value = -24364.5
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H86").Value = 397.66666
$ws.Range("J86").Value = 397.66666
$ws.Range("L86").Value = 1192.99998
$ws.Range("N86").Value = -3564.99998
$ws.Range("H88").Value = 29250
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 29250
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 87750
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -88606
$ws.Range("H89").Value = 397.66666
$ws.Range("J89").Value = 397.66666
$ws.Range("L89").Value = 3578.99994
$ws.Range("N89").Value = -15434.99994
$ws.Range("H91").Value = 29250
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 29250
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 87750
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -90714
$ws.Range("H101").Value = 12511312
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 12511312
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 37533936
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -37538804
$ws.Range("H121").Value = 3596.7646
$ws.Range("I121").Value = 444.66666
$ws.Range("J121").Value = 5316.091
$ws.Range("K121").Value = 1333.99998
$ws.Range("L121").Value = 15948.273
$ws.Range("M121").Value = -23.99998000000005
$ws.Range("N121").Value = -18568.273
$ws.Range("H131").Value = 2247.6365
$ws.Range("J131").Value = 2247.6365
$ws.Range("L131").Value = 6742.9095
$ws.Range("N131").Value = -16822.9095

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 9999.5
$ws.Range("I18").Value = 9999
$ws.Range("K18").Value = 9999
$ws.Range("M18").Value = -9706
$ws.Range("H21").Value = 8001
$ws.Range("I21").Value = 8001
$ws.Range("K21").Value = 8001
$ws.Range("M21").Value = -7828
$ws.Range("H24").Value = 5035.4287
$ws.Range("J24").Value = 5035.4287
$ws.Range("L24").Value = 5035.4287
$ws.Range("N24").Value = -5381.4287
$ws.Range("H30").Value = 8001
$ws.Range("I30").Value = 8001
$ws.Range("K30").Value = 8001
$ws.Range("M30").Value = -7896
$ws.Range("H70").Value = 61426.8
$ws.Range("I70").Value = 61103.707
$ws.Range("J70").Value = 61731.945
$ws.Range("K70").Value = 61103.707
$ws.Range("L70").Value = 61731.945
$ws.Range("M70").Value = -60833.707
$ws.Range("N70").Value = -62271.945
$ws.Range("H73").Value = 61426.8
$ws.Range("I73").Value = 61103.707
$ws.Range("J73").Value = 61731.945
$ws.Range("K73").Value = 61103.707
$ws.Range("L73").Value = 61731.945
$ws.Range("M73").Value = -60167.707
$ws.Range("N73").Value = -63603.945
$ws.Range("H101").Value = 55017
$ws.Range("J101").Value = 55017
$ws.Range("L101").Value = 55017
$ws.Range("N101").Value = -61507
$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524
$ws.Range("H122").Value = 51403.855
$ws.Range("I122").Value = 80652.54
$ws.Range("J122").Value = 3874.75
$ws.Range("K122").Value = 241957.62
$ws.Range("L122").Value = 11624.25
$ws.Range("M122").Value = -239507.62
$ws.Range("N122").Value = -16524.25
$ws.Range("H126").Value = 2321.842
$ws.Range("I126").Value = 2192.8572
$ws.Range("K126").Value = 6578.571599999999
$ws.Range("M126").Value = -4108.571599999999
$ws.Range("H132").Value = 11544.541
$ws.Range("I132").Value = 12826.5625
$ws.Range("J132").Value = 3339.6
$ws.Range("K132").Value = 38479.6875
$ws.Range("L132").Value = 10018.8
$ws.Range("M132").Value = -35949.6875
$ws.Range("N132").Value = -15078.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3686.8235
$ws.Range("I7").Value = 3126
$ws.Range("K7").Value = 3126
$ws.Range("M7").Value = -3014
$ws.Range("H22").Value = 4032.6667
$ws.Range("I22").Value = 3401
$ws.Range("J22").Value = 4090.0908
$ws.Range("K22").Value = 3401
$ws.Range("L22").Value = 4090.0908
$ws.Range("M22").Value = -3106
$ws.Range("N22").Value = -4680.0908
$ws.Range("H27").Value = 4032.6667
$ws.Range("I27").Value = 3401
$ws.Range("J27").Value = 4090.0908
$ws.Range("K27").Value = 3401
$ws.Range("L27").Value = 4090.0908
$ws.Range("M27").Value = -3294
$ws.Range("N27").Value = -4304.0908
$ws.Range("H40").Value = 3229.6155
$ws.Range("I40").Value = 1982.5
$ws.Range("K40").Value = 1982.5
$ws.Range("M40").Value = -1846.5
$ws.Range("H68").Value = 1684
$ws.Range("J68").Value = 1684
$ws.Range("L68").Value = 1684
$ws.Range("N68").Value = -3182
$ws.Range("H71").Value = 1684
$ws.Range("J71").Value = 1684
$ws.Range("L71").Value = 8420
$ws.Range("N71").Value = -15908
$ws.Range("H100").Value = 2280.5
$ws.Range("I100").Value = 1781.6
$ws.Range("K100").Value = 1781.6
$ws.Range("M100").Value = -1240.6
$ws.Range("H126").Value = 3686.8235
$ws.Range("I126").Value = 3126
$ws.Range("K126").Value = 9378
$ws.Range("M126").Value = -6908
$ws.Range("H136").Value = 7917.2915
$ws.Range("I136").Value = 4316.5
$ws.Range("K136").Value = 12949.5
$ws.Range("M136").Value = -10399.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17075.5
$ws.Range("J45").Value = 17075
$ws.Range("L45").Value = 17075
$ws.Range("N45").Value = -18057
$ws.Range("H62").Value = 3085.25
$ws.Range("I62").Value = 3085.25
$ws.Range("K62").Value = 3085.25
$ws.Range("M62").Value = -2461.25
$ws.Range("H65").Value = 3085.25
$ws.Range("I65").Value = 3085.25
$ws.Range("K65").Value = 15426.25
$ws.Range("M65").Value = -12306.25
$ws.Range("H107").Value = 859.5122
$ws.Range("I107").Value = 807.1429
$ws.Range("J107").Value = 1165
$ws.Range("K107").Value = 2421.4287
$ws.Range("L107").Value = 3495
$ws.Range("M107").Value = -501.4287000000004
$ws.Range("N107").Value = -7335
$ws.Range("H122").Value = 68411.82
$ws.Range("I122").Value = 2328.7
$ws.Range("K122").Value = 6986.099999999999
$ws.Range("M122").Value = -4536.099999999999

